$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 13.86384647080068, 0, 17.50910633199374)
    3 = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1, 1.104883657715537)
    4 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
    5 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 0, 21.98653043760045)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
